# Update PLC data 2025-10-13 14:12:31
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 286
$ws.Range("C3").Value = 176265
$ws.Range("C4").Value = 166221
$ws.Range("C8").Value = 64.58
